# Mexico GP comparison table: re-ran the notebook and picked up the fix to
# the overtaking-count bug, which changed several rows' computed values and
# re-sorted a few rows by position_sim (Oscar Piastri/George Russell and
# Franco Colapinto/Esteban Ocon/Valtteri Bottas/Liam Lawson swap places).
# Columns: 1=driver_name 2=laps_completed 3=position_sim 4=position_actual
#          5=overtake_error 6=cumulative_time_error 7=gap_error
# Row 1 is the header row, so data rows are table rows 2-21.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2, 5).Range.Text = "6"
$t.Cell(2, 6).Range.Text = "239.362"
$t.Cell(3, 6).Range.Text = "237.310"
$t.Cell(3, 7).Range.Text = "-2.052"
$t.Cell(4, 5).Range.Text = "6"
$t.Cell(4, 6).Range.Text = "230.995"
$t.Cell(4, 7).Range.Text = "-8.367"
$t.Cell(5, 6).Range.Text = "231.999"
$t.Cell(5, 7).Range.Text = "-7.363"
$t.Cell(6, 1).Range.Text = "Oscar Piastri"
$t.Cell(6, 4).Range.Text = "8"
$t.Cell(6, 6).Range.Text = "226.703"
$t.Cell(6, 7).Range.Text = "-12.659"
$t.Cell(7, 1).Range.Text = "George Russell"
$t.Cell(7, 4).Range.Text = "5"
$t.Cell(7, 6).Range.Text = "240.531"
$t.Cell(7, 7).Range.Text = "1.169"
$t.Cell(8, 5).Range.Text = "-2"
$t.Cell(8, 6).Range.Text = "233.959"
$t.Cell(8, 7).Range.Text = "-5.403"
$t.Cell(9, 5).Range.Text = "6"
$t.Cell(9, 6).Range.Text = "241.197"
$t.Cell(9, 7).Range.Text = "1.835"
$t.Cell(10, 6).Range.Text = "319.189"
$t.Cell(10, 7).Range.Text = "-0.758"
$t.Cell(11, 5).Range.Text = "-1"
$t.Cell(11, 6).Range.Text = "403.960"
$t.Cell(11, 7).Range.Text = "84.013"
$t.Cell(12, 5).Range.Text = "-3"
$t.Cell(12, 6).Range.Text = "323.260"
$t.Cell(12, 7).Range.Text = "3.313"
$t.Cell(13, 1).Range.Text = "Franco Colapinto"
$t.Cell(13, 4).Range.Text = "12"
$t.Cell(13, 5).Range.Text = "-2"
$t.Cell(13, 6).Range.Text = "322.128"
$t.Cell(13, 7).Range.Text = "2.181"
$t.Cell(14, 1).Range.Text = "Esteban Ocon"
$t.Cell(14, 4).Range.Text = "13"
$t.Cell(14, 5).Range.Text = "7"
$t.Cell(14, 6).Range.Text = "317.590"
$t.Cell(14, 7).Range.Text = "-2.357"
$t.Cell(15, 1).Range.Text = "Valtteri Bottas"
$t.Cell(15, 4).Range.Text = "14"
$t.Cell(15, 5).Range.Text = "-4"
$t.Cell(15, 6).Range.Text = "328.555"
$t.Cell(15, 7).Range.Text = "8.608"
$t.Cell(16, 5).Range.Text = "1"
$t.Cell(16, 6).Range.Text = "309.552"
$t.Cell(16, 7).Range.Text = "-10.395"
$t.Cell(17, 1).Range.Text = "Liam Lawson"
$t.Cell(17, 4).Range.Text = "16"
$t.Cell(17, 5).Range.Text = "-5"
$t.Cell(17, 6).Range.Text = "307.100"
$t.Cell(17, 7).Range.Text = "-12.847"
$t.Cell(18, 5).Range.Text = "4"
$t.Cell(18, 6).Range.Text = "317.930"
$t.Cell(18, 7).Range.Text = "-2.017"
